$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46000.94930555556
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B2").Value = -8.01
$ws.Range("C2").Value = 114.18
$ws.Range("D2").Value = 2.9
$ws.Range("E2").Value = 17
$ws.Range("G2").Value = 0.9450121840697095
$ws.Range("H2").Value = 94.5
$ws.Range("I2").Value = "Terdampak"
$ws.Range("J2").Value = 5.877254504294055

$ws.Range("A3").Value = 46000.94930555556
$ws.Range("A3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B3").Value = -8.01
$ws.Range("C3").Value = 114.18
$ws.Range("D3").Value = 2.9
$ws.Range("E3").Value = 17
$ws.Range("G3").Value = 0.9511983133618672
$ws.Range("H3").Value = 95.12
$ws.Range("I3").Value = "Terdampak"
$ws.Range("J3").Value = 5.88959993921047

$ws.Range("A4").Value = 46000.96944444445
$ws.Range("A4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B4").Value = -7.96
$ws.Range("C4").Value = 112.77
$ws.Range("D4").Value = 2.1
$ws.Range("E4").Value = 8
$ws.Range("G4").Value = 0.0001
$ws.Range("H4").Value = 0.01
$ws.Range("I4").Value = "Aman"
$ws.Range("J4").Value = 3

$ws.Range("A5").Value = 46000.96944444445
$ws.Range("A5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B5").Value = -7.96
$ws.Range("C5").Value = 112.77
$ws.Range("D5").Value = 2.1
$ws.Range("E5").Value = 8
$ws.Range("G5").Value = 0.05978295139323081
$ws.Range("H5").Value = 5.98
$ws.Range("I5").Value = "Aman"
$ws.Range("J5").Value = 3

$ws.Range("A6").Value = 46002.14027777778
$ws.Range("A6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B6").Value = -7.98
$ws.Range("C6").Value = 114.15
$ws.Range("D6").Value = 2.2
$ws.Range("E6").Value = 4
$ws.Range("G6").Value = 0.8662135250639107
$ws.Range("H6").Value = 86.62
$ws.Range("I6").Value = "Terdampak"
$ws.Range("J6").Value = 3.99418365778594

$ws.Range("A7").Value = 46002.14027777778
$ws.Range("A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B7").Value = -7.98
$ws.Range("C7").Value = 114.15
$ws.Range("D7").Value = 2.2
$ws.Range("E7").Value = 4
$ws.Range("G7").Value = 0.7810637733539998
$ws.Range("H7").Value = 78.11
$ws.Range("I7").Value = "Terdampak"
$ws.Range("J7").Value = 3.875524059060939

$ws.Range("A8").Value = 46002.91805555556
$ws.Range("A8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B8").Value = -7.92
$ws.Range("C8").Value = 114.16
$ws.Range("D8").Value = 2.2
$ws.Range("E8").Value = 6
$ws.Range("G8").Value = 0.9175182761045644
$ws.Range("H8").Value = 91.75
$ws.Range("I8").Value = "Terdampak"
$ws.Range("J8").Value = 4.065678888822081

$ws.Range("A9").Value = 46002.91805555556
$ws.Range("A9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B9").Value = -7.92
$ws.Range("C9").Value = 114.16
$ws.Range("D9").Value = 2.2
$ws.Range("E9").Value = 6
$ws.Range("G9").Value = 0.8796969689600115
$ws.Range("H9").Value = 87.97
$ws.Range("I9").Value = "Terdampak"
$ws.Range("J9").Value = 4.012973378366607

$ws.Range("A10").Value = 46003.22152777778
$ws.Range("A10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B10").Value = -7.98
$ws.Range("C10").Value = 114.18
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 13
$ws.Range("G10").Value = 0.7823351233528767
$ws.Range("H10").Value = 78.23
$ws.Range("I10").Value = "Terdampak"
$ws.Range("J10").Value = 5.802792935549129

$ws.Range("A11").Value = 46003.22152777778
$ws.Range("A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B11").Value = -7.98
$ws.Range("C11").Value = 114.18
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 13
$ws.Range("G11").Value = 0.8780029949326082
$ws.Range("H11").Value = 87.8
$ws.Range("I11").Value = "Terdampak"
$ws.Range("J11").Value = 6.002316294436578

$ws.Range("A12").Value = 46003.23680555556
$ws.Range("A12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B12").Value = -7.97
$ws.Range("C12").Value = 114.13
$ws.Range("D12").Value = 2.2
$ws.Range("E12").Value = 6
$ws.Range("G12").Value = 0.911968362286605
$ws.Range("H12").Value = 91.2
$ws.Range("I12").Value = "Terdampak"
$ws.Range("J12").Value = 4.057944861025488

$ws.Range("A13").Value = 46003.76111111111
$ws.Range("A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B13").Value = -7.91
$ws.Range("C13").Value = 114.14
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 14
$ws.Range("G13").Value = 0.5636424251211228
$ws.Range("H13").Value = 56.36
$ws.Range("I13").Value = "Aman"
$ws.Range("J13").Value = 3.156214049501883

$ws.Range("A14").Value = 46004.62708333333
$ws.Range("A14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B14").Value = -7.98
$ws.Range("C14").Value = 114.24
$ws.Range("D14").Value = 3.1
$ws.Range("E14").Value = 10
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 100
$ws.Range("I14").Value = "Terdampak"
$ws.Range("J14").Value = 6.529222387025386

$ws.Range("A15").Value = 46004.62986111111
$ws.Range("A15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B15").Value = -7.99
$ws.Range("C15").Value = 114.21
$ws.Range("D15").Value = 2.9
$ws.Range("E15").Value = 5
$ws.Range("G15").Value = 0.5371639818466017
$ws.Range("H15").Value = 53.72
$ws.Range("I15").Value = "Aman"
$ws.Range("J15").Value = 5.063326569954734

$ws.Range("A16").Value = 46004.62986111111
$ws.Range("A16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B16").Value = -7.99
$ws.Range("C16").Value = 114.21
$ws.Range("D16").Value = 2.9
$ws.Range("E16").Value = 5
$ws.Range("G16").Value = 0.6581567157044004
$ws.Range("H16").Value = 65.82
$ws.Range("I16").Value = "Aman"
$ws.Range("J16").Value = 5.304787401220269

$ws.Range("A17").Value = 46004.97708333333
$ws.Range("A17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B17").Value = -7.99
$ws.Range("C17").Value = 114.22
$ws.Range("D17").Value = 2.9
$ws.Range("E17").Value = 5
$ws.Range("G17").Value = 0.5146274993287281
$ws.Range("H17").Value = 51.46
$ws.Range("I17").Value = "Aman"
$ws.Range("J17").Value = 5.018351325334156

$ws.Range("A18").Value = 46004.97708333333
$ws.Range("A18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B18").Value = -7.99
$ws.Range("C18").Value = 114.22
$ws.Range("D18").Value = 2.9
$ws.Range("E18").Value = 5
$ws.Range("G18").Value = 0.4953046286158333
$ws.Range("H18").Value = 49.53
$ws.Range("I18").Value = "Aman"
$ws.Range("J18").Value = 4.979789369785192

$ws.Range("A19").Value = 46004.99583333333
$ws.Range("A19").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B19").Value = -7.97
$ws.Range("C19").Value = 114.16
$ws.Range("D19").Value = 2.3
$ws.Range("E19").Value = 6
$ws.Range("G19").Value = 0.7863398500639497
$ws.Range("H19").Value = 78.63
$ws.Range("I19").Value = "Terdampak"
$ws.Range("J19").Value = 4.113867288164345

$ws.Range("A20").Value = 46006.00833333333
$ws.Range("A20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B20").Value = -7.97
$ws.Range("C20").Value = 114.15
$ws.Range("D20").Value = 2.2
$ws.Range("E20").Value = 11
$ws.Range("G20").Value = 0.8449057463909231
$ws.Range("H20").Value = 84.49
$ws.Range("I20").Value = "Terdampak"
$ws.Range("J20").Value = 3.964490412458913

$ws.Range("A21").Value = 46006.48263888889
$ws.Range("A21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B21").Value = -7.99
$ws.Range("C21").Value = 114.19
$ws.Range("D21").Value = 2.6
$ws.Range("E21").Value = 12
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 100
$ws.Range("I21").Value = "Terdampak"
$ws.Range("J21").Value = 5.1946538276055

$ws.Range("A22").Value = 46006.86041666667
$ws.Range("A22").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B22").Value = -7.96
$ws.Range("C22").Value = 114.15
$ws.Range("D22").Value = 2.5
$ws.Range("E22").Value = 10
$ws.Range("G22").Value = 0.8145592135071086
$ws.Range("H22").Value = 81.46
$ws.Range("I22").Value = "Terdampak"
$ws.Range("J22").Value = 4.631294577132157
